$wb = $excel.ActiveWorkbook
$ft = $wb.Worksheets.Item("FT")
$ft.Activate()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 4
$ft.Range("B66").Select()
Write-Host "done"
